$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update username/password values for rows 2-6 to the new manager account
$ws.Range("A2").Value = "mngr429183"
$ws.Range("B2").Value = "Avamyze"
$ws.Range("A3").Value = "mngr429183"
$ws.Range("B3").Value = "Avamyze"
$ws.Range("A4").Value = "mngr429183"
$ws.Range("B4").Value = "Avamyze"
$ws.Range("A5").Value = "mngr429183"
$ws.Range("B5").Value = "Avamyze"
$ws.Range("A6").Value = "mngr429183"
$ws.Range("B6").Value = "Avamyze"

# Update selection to A8
$ws.Range("A8").Select()
